$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value (all target cells store text, so we force
# a "@" text number-format before assigning, then reset the style to
# "Normal" so no stray per-cell number-format sticks around afterwards.
$updates = [ordered]@{
    "D2" = "297.94"
    "E2" = "-1.97%"
    "G2" = "2"
    "D3" = "31.38"
    "E3" = "-1.19%"
    "G3" = "2"
    "D4" = "5.098"
    "E4" = "-2.09%"
    "G4" = "2"
    "D5" = "0.07933"
    "E5" = "6.72%"
    "G5" = "2"
    "D6" = "2.315"
    "E6" = "3.80%"
    "G6" = "2"
    "D7" = "7.739"
    "E7" = "-2.75%"
    "G7" = "2"
    "D8" = "3.886"
    "E8" = "0.57%"
    "G8" = "2"
    "D9" = "0.9240"
    "E9" = "0.67%"
    "G9" = "2"
    "D10" = "0.1736"
    "E10" = "0.27%"
    "G10" = "2"
    "D11" = "0.07401"
    "E11" = "-3.57%"
    "G11" = "2"
    "D12" = "0.08948"
    "E12" = "10.02%"
    "G12" = "2"
    "D13" = "0.03017"
    "E13" = "1.26%"
    "G13" = "2"
    "D14" = "0.1000"
    "E14" = "0.83%"
    "G14" = "2"
    "D15" = "0.001516"
    "E15" = "1.72%"
    "G15" = "2"
    "D16" = "0.006099"
    "E16" = "0.11%"
    "G16" = "2"
    "D17" = "3.498"
    "E17" = "0.24%"
    "G17" = "2"
    "D18" = "2.288"
    "E18" = "2.67%"
    "G18" = "2"
    "E19" = "0.26%"
    "G19" = "2"
    "D20" = "0.1329"
    "E20" = "-1.14%"
    "G20" = "2"
    "D21" = "4.155"
    "E21" = "-10.40%"
    "G21" = "2"
    "D22" = "0.1693"
    "E22" = "8.57%"
    "G22" = "2"
    "D23" = "0.04610"
    "E23" = "-0.29%"
    "G23" = "2"
    "D24" = "0.001242"
    "E24" = "1.66%"
    "G24" = "2"
    "D25" = "0.004446"
    "E25" = "-1.18%"
    "G25" = "2"
    "E26" = "-7.25%"
    "G26" = "2"
    "D27" = "0.0003400"
    "E27" = "24.49%"
    "G27" = "2"
    "G28" = "2"
    "G29" = "2"
    "G30" = "2"
    "G31" = "2"
    "G32" = "2"
    "G33" = "2"
    "G34" = "2"
    "G35" = "2"
    "G36" = "2"
    "G37" = "2"
    "G38" = "2"
    "D39" = "0.01741"
    "E39" = "-2.79%"
    "G39" = "2"
    "E40" = "0.79%"
    "G40" = "2"
    "D41" = "0.006962"
    "E41" = "-5.13%"
    "G41" = "2"
    "D42" = "0.1355"
    "E42" = "0.04%"
    "G42" = "2"
    "D43" = "0.002191"
    "E43" = "1.40%"
    "G43" = "2"
    "D44" = "0.01026"
    "E44" = "-6.14%"
    "G44" = "2"
    "D45" = "0.00006319"
    "E45" = "1.31%"
    "G45" = "2"
    "D46" = "0.00000000750"
    "E46" = "0.05%"
    "G46" = "2"
    "D47" = "0.007985"
    "E47" = "-19.01%"
    "G47" = "2"
    "D48" = "0.7478"
    "E48" = "-8.87%"
    "G48" = "2"
    "D49" = "0.00002101"
    "E49" = "0.05%"
    "G49" = "2"
    "D50" = "0.0002001"
    "E50" = "0.12%"
    "G50" = "2"
    "G51" = "2"
}

foreach ($cellRef in $updates.Keys) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cellRef]
    $rng.Style = "Normal"
}
